$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: add the two trailing columns that were missing ---
$ws.Range("X3").Value = -0.34999899999999684
$ws.Range("Y3").Value = "Down"

# --- Row 4: new trade row ---
$ws.Range("A4").Value = 42649.612187500003
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = "Buy"
$ws.Range("D4").Value = 36
$ws.Range("E4").Value = 7436
$ws.Range("F4").Value = 787
$ws.Range("G4").Value = 63
$ws.Range("H4").Value = 35
$ws.Range("I4").Value = 86
$ws.Range("J4").Value = 13
$ws.Range("K4").Value = 5352
$ws.Range("L4").Value = 110
$ws.Range("M4").Value = 61
$ws.Range("N4").Value = 19
$ws.Range("O4").Value = 3
$ws.Range("P4").Value = "Bag"
$ws.Range("Q4").Value = 35.483823948801813
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0.0965
$ws.Range("S4").NumberFormat = "0.00%"
$ws.Range("T4").Value = 0.0269
$ws.Range("T4").NumberFormat = "0.00%"
$ws.Range("U4").Value = 4.82
$ws.Range("V4").Value = 2.2799999999999998
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = -0.34999899999999684
$ws.Range("Y4").Value = "Down"

# --- Row 5: new trade row ---
$ws.Range("A5").Value = 42649.635567129626
$ws.Range("A5").NumberFormat = "m/d/yy h:mm"
$ws.Range("B5").Value = 11
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 46
$ws.Range("E5").Value = 9051
$ws.Range("F5").Value = 1047
$ws.Range("G5").Value = 66
$ws.Range("H5").Value = 32
$ws.Range("I5").Value = 90
$ws.Range("J5").Value = 9
$ws.Range("K5").Value = 6478
$ws.Range("L5").Value = 150
$ws.Range("M5").Value = 74
$ws.Range("N5").Value = 69
$ws.Range("O5").Value = 7
$ws.Range("P5").Value = "Bag"
$ws.Range("Q5").Value = 35.483823948801813
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0.0965
$ws.Range("S5").NumberFormat = "0.00%"
$ws.Range("T5").Value = 0.0269
$ws.Range("T5").NumberFormat = "0.00%"
$ws.Range("U5").Value = 4.82
$ws.Range("V5").Value = 2.2799999999999998
$ws.Range("W5").Value = 0

Write-Host "applied"
